# reference.docx: use current standard Word theme.
# Switches the document theme's color scheme and major/minor Latin fonts
# from the legacy Office 2007-2010 "Office" theme (Calibri/Cambria,
# blue-ish accents) to the current default Word theme (Aptos
# Display/Aptos, the modern accent palette), and drops the stray unused
# "Default Paragraph Font" latent style definition left over from the
# old reference template.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Theme color scheme: dk2/lt2/accent1-6/hlink/folHlink get the new
#    palette (dk1/lt1 - windowText/window - are unchanged).
# ---------------------------------------------------------------------

function HexToRGB($hex) {
    # Word's ColorFormat.RGB / ThemeColor.RGB is 0x00BBGGRR (little-endian),
    # same as the VBA RGB() function.
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$colorScheme = $d.DocumentTheme.ThemeColorScheme

# Index order follows MsoThemeColorSchemeIndex:
#  1 dk1  2 lt1  3 dk2  4 lt2  5-10 accent1-6  11 hlink  12 folHlink
$colorScheme.Colors(3).RGB = HexToRGB "0E2841"   # dk2       1F497D -> 0E2841
$colorScheme.Colors(4).RGB = HexToRGB "E8E8E8"   # lt2       EEECE1 -> E8E8E8
$colorScheme.Colors(5).RGB = HexToRGB "156082"   # accent1   4F81BD -> 156082
$colorScheme.Colors(6).RGB = HexToRGB "E97132"   # accent2   C0504D -> E97132
$colorScheme.Colors(7).RGB = HexToRGB "196B24"   # accent3   9BBB59 -> 196B24
$colorScheme.Colors(8).RGB = HexToRGB "0F9ED5"   # accent4   8064A2 -> 0F9ED5
$colorScheme.Colors(9).RGB = HexToRGB "A02B93"   # accent5   4BACC6 -> A02B93
$colorScheme.Colors(10).RGB = HexToRGB "4EA72E"  # accent6   F79646 -> 4EA72E
$colorScheme.Colors(11).RGB = HexToRGB "467886"  # hlink     0000FF -> 467886
$colorScheme.Colors(12).RGB = HexToRGB "96607D"  # folHlink  800080 -> 96607D

# ---------------------------------------------------------------------
# 2. Theme font scheme: swap the serif "Cambria"/"Calibri" pairing for
#    the new "Aptos Display"/"Aptos" pairing used by the current
#    default Word theme.
# ---------------------------------------------------------------------

$fontScheme = $d.DocumentTheme.ThemeFontScheme

$fontScheme.MajorFont.Latin = "Aptos Display"
$fontScheme.MajorFont.Panose = "02110004020202020204"

$fontScheme.MinorFont.Latin = "Aptos"
$fontScheme.MinorFont.Panose = "02110004020202020204"

# ---------------------------------------------------------------------
# 3. Drop the unused built-in "Default Paragraph Font" latent style
#    definition (still referenced elsewhere as a basedOn target, but no
#    longer declared itself).
# ---------------------------------------------------------------------

$d.Styles("Default Paragraph Font").Delete()
